# Updates odds/value cells on Sheet1 to reflect the latest FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Bristol City vs Sheffield Utd)
$ws.Range("G2").Value = 3.3
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 1.83
$ws.Range("AA2").Value = 29
$ws.Range("AG2").Value = 301
$ws.Range("AJ2").Value = 9.5
$ws.Range("AX2").Value = 13

# Row 5 (Sheffield Wed vs Norwich)
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.85

# Row 6 (Plymouth vs Portsmouth)
$ws.Range("G6").Value = 2.3
$ws.Range("H6").Value = 3.6
$ws.Range("I6").Value = 2.88
$ws.Range("J6").Value = 2.88
$ws.Range("K6").Value = 2.3
$ws.Range("Q6").Value = 1.68
$ws.Range("R6").Value = 2.05
$ws.Range("U6").Value = 1.57
$ws.Range("V6").Value = 2.25
$ws.Range("AO6").Value = 12
$ws.Range("AP6").Value = 19

# Row 7 (Swansea vs Watford)
$ws.Range("Q7").Value = 1.62
$ws.Range("R7").Value = 2.15

# Row 8 (Burton vs Crawley)
$ws.Range("Q8").Value = 1.77

# Row 9 (Stockport County vs Wycombe)
$ws.Range("Q9").Value = 1.95
$ws.Range("R9").Value = 1.9
$ws.Range("S9").Value = 1.4

# Row 10 (Sakhnin vs Beitar Jerusalem)
$ws.Range("G10").Value = 4.1
$ws.Range("I10").Value = 1.67
$ws.Range("J10").Value = 4.5
$ws.Range("L10").Value = 2.25
$ws.Range("N10").Value = 17
$ws.Range("Q10").Value = 1.6
$ws.Range("R10").Value = 2.3
$ws.Range("S10").Value = 1.3
$ws.Range("U10").Value = 1.62
$ws.Range("V10").Value = 2.2
$ws.Range("X10").Value = 23
$ws.Range("Y10").Value = 13
$ws.Range("Z10").Value = 41
$ws.Range("AA10").Value = 29
$ws.Range("AI10").Value = 9.5
$ws.Range("AL10").Value = 13
